$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195112705230713
$ws.Range("B1").Value = 1.673968195915222
$ws.Range("C1").Value = 6.764725208282471
$ws.Range("D1").Value = 2.254708290100098
$ws.Range("E1").Value = 1.186519503593445
